$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.321.28'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '3.012.22'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'355.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").Value = "'108.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.80%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D9").Value = "'0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").Value = "'38.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.70%  '
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("E12").Value = '  -4.50%  '
$ws.Range("D13").Value = "'19.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("D14").Value = '3.487.45'
$ws.Range("E14").Value = '  +1.93%  '
$ws.Range("D15").Value = "'7.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").Value = '3.004.18'
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").Value = '52.364.63'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = "'3.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.26%  '
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("D21").Value = "'13.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.34%  '
$ws.Range("D22").Value = '0.0₃0977'
$ws.Range("E22").Value = '  -1.53%  '
$ws.Range("D23").Value = "'69.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").Value = "'265.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.67%  '
$ws.Range("D25").Value = "'2.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").Value = "'0.179"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").Value = "'27.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("D28").Value = "'7.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("E30").Value = '  -3.11%  '
$ws.Range("D31").Value = "'6.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("D32").Value = "'10.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.67%  '
$ws.Range("D33").Value = "'36.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("E34").Value = '  +19.04%  '
$ws.Range("D35").Value = "'50.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.87%  '
$ws.Range("D36").Value = "'0.0443"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = "'3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("D39").Value = "'2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("D40").Value = "'18.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.53%  '
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("D43").Value = "'23.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("D44").Value = "'123.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.14%  '
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("D46").Value = '2.135.57'
$ws.Range("E46").Value = '  -1.36%  '
$ws.Range("D47").Value = "'3.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -5.28%  '
$ws.Range("D49").Value = '3.313.14'
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("D51").Value = "'0.0333"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.42%  '
